$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Espinaca" table. It belongs
# chronologically before the existing row 63, so insert a fresh row at 63 -
# this shifts the former rows 63:88 down to 64:89 (matching the diff) and
# keeps the inherited column-D date style (s="2") on the new row.
$ws.Rows(63).Insert()

$row = 63
$ws.Cells.Item($row, 1).Value  = 1
$ws.Cells.Item($row, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value  = 44917
$ws.Cells.Item($row, 5).Value  = 15
$ws.Cells.Item($row, 6).Value  = 100112012
$ws.Cells.Item($row, 7).Value  = "Espinaca"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 450
$ws.Cells.Item($row, 11).Value = 4000
$ws.Cells.Item($row, 12).Value = 4500
$ws.Cells.Item($row, 13).Value = 4222
$ws.Cells.Item($row, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 1407
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
